# Update the "Fitness" values (column C) for rows 2 through 69 on the active sheet.
# New values correspond to the updated run_5.xlsx data described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    12843,11397,10806,10806,10806,10806,10806,9925,9925,9925,
    9925,8583,8583,8583,8583,8583,8583,8583,8497,8497,
    8497,8497,8497,8422,8422,7910,7910,7910,7573,7573,
    7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,
    7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,
    7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,
    7573,7573,7573,7573,7573,7573,7573,7573
)

$startRow = 2
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 3).Value = $newValues[$i]
}
